# Generate Report for Handback
#
# Replaces the two handback file identifiers (UUID-based markdown file
# names) and their derived xliff / timestamp values across the three
# worksheets (Overview, zh-cn, de-de), matching a newer CI run.

$wb = $excel.ActiveWorkbook

# ---- old -> new identifiers -------------------------------------------------
# (old) 5950f920-e0d3-4c1e-9147-b86f18c5a9cf.md -> (new) 93b4c6ab-5540-4dd3-b33c-fc9a648f8162.md
# (old) dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51.md -> (new) ffffe88c0b38-6372-4f48-8d62-7a81284afec1.md
$newGuid1 = "93b4c6ab-5540-4dd3-b33c-fc9a648f8162"
$newGuid2 = "ffffe88c0b38-6372-4f48-8d62-7a81284afec1"

$newHash  = "f9615f04b769825cba53aa39339543ff522f2111"

$newGenDate   = "2016-08-27 06:59:19"   # Overview col G / zh-cn+de-de col H
$newZhStart   = "2016-08-27 06:59:14"   # zh-cn col H (correspond handoff datetime)
$newZhBack    = "2016-08-27 06:59:31"   # zh-cn col K (correspond handback datetime)
$newDeBack    = "2016-08-27 06:59:37"   # de-de col K (correspond handback datetime)

$newZhXlf = "$newGuid1.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid1.$newHash.de-de.xlf"

# =========================================================================
# Overview sheet
# =========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = "e2e\$newGuid1.md"
$ws.Range("G2").Value = $newGenDate

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "e2e\$newGuid2.md"
$ws.Range("G3").Value = $newGenDate

# Rebuild the hyperlinks on B2/B3 so their display text matches the new
# file names (underlying targets are left untouched).
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid1.md", "", "", "e2e\$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md")

# =========================================================================
# zh-cn sheet
# =========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = $newZhXlf
$ws.Range("H2").Value = $newZhStart
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = $newZhXlf
$ws.Range("K2").Value = $newZhBack

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = $newZhXlf
$ws.Range("H3").Value = $newZhStart
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = $newZhXlf
$ws.Range("K3").Value = $newZhBack

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a9496254b3f34e3eca08c4678f70cdc3ef300bc/e2e/$newGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a9496254b3f34e3eca08c4678f70cdc3ef300bc/e2e/$newGuid2.md", "", "", "$newGuid2.md")

# =========================================================================
# de-de sheet
# =========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = $newDeXlf
$ws.Range("H2").Value = $newGenDate
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = $newDeXlf
$ws.Range("K2").Value = $newDeBack

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = $newDeXlf
$ws.Range("H3").Value = $newGenDate
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = $newDeXlf
$ws.Range("K3").Value = $newDeBack

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f1c3ff5529e272061bb14a9d8997e7a3b5b20ca5/e2e/$newGuid1.md", "", "", "$newGuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$newGuid2.md", "", "", "$newGuid2.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f1c3ff5529e272061bb14a9d8997e7a3b5b20ca5/e2e/$newGuid2.md", "", "", "$newGuid2.md")
